$d = $word.ActiveDocument

# --- 1) First paragraph: "=Critère A" -> "Critère A" (drop the leading "=" run) ---
$p1 = $d.Paragraphs.Item(1)
$xml1 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:pPr><w:jc w:val='both'/><w:rPr><w:lang w:val='fr-CA'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='fr-CA'/></w:rPr><w:t>Critère A</w:t></w:r>" +
        "</w:p>"
$p1.Range.InsertXML($xml1)

# --- 2) "Aspect ii : " paragraph: wrap "Aspect ii" in proofErr gramStart/gramEnd,
#        split off the " : " (with leading non-breaking space) into its own run ---
$p2 = $d.Paragraphs.Item(5)
$xml2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:pPr><w:jc w:val='both'/><w:rPr><w:lang w:val='fr-CA'/></w:rPr></w:pPr>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:rPr><w:lang w:val='fr-CA'/></w:rPr><w:t>Aspect ii</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:rPr><w:lang w:val='fr-CA'/></w:rPr><w:t xml:space='preserve'>&#160;: </w:t></w:r>" +
        "</w:p>"
$p2.Range.InsertXML($xml2)

# --- 3) "Aspect iii :" paragraph: same treatment, no trailing space after colon ---
$p3 = $d.Paragraphs.Item(8)
$xml3 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:pPr><w:jc w:val='both'/><w:rPr><w:lang w:val='fr-CA'/></w:rPr></w:pPr>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:rPr><w:lang w:val='fr-CA'/></w:rPr><w:t>Aspect iii</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:rPr><w:lang w:val='fr-CA'/></w:rPr><w:t>&#160;:</w:t></w:r>" +
        "</w:p>"
$p3.Range.InsertXML($xml3)
